# "Fixed shell bug, removed duplicate C12 line"
#
# The sheet had a duplicated "C12" row (row 12) — a leftover copy of the
# C12 entry already present at row 4, complete with a stray note in I12
# ("C12 density assumes in graphite form"). Remove that duplicate row so
# every row below it shifts up by one; Excel automatically re-indexes the
# shared formulas (E3:E14 -> E3:E13), drops the now-unreferenced shared
# string, and renumbers the remaining string indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate "C12" row (row 12).
$ws.Rows.Item(12).Delete() | Out-Null

# Leave the selection where the author left it after the edit.
$ws.Range("C11").Select() | Out-Null
